$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fields on row 2 (plain string replacements) ---
$ws.Range("A2").Value = "merio"
$ws.Range("B2").Value = "aloni"
$ws.Range("D2").Value = "merio@test.com"
$ws.Range("H2").Value = "c com"
$ws.Range("L2").Value = "B+"
$ws.Range("M2").Value = "nehru park"
$ws.Range("O2").Value = "rajasthan"
$ws.Range("P2").Value = "pali"

# --- Date-look-alike text fields: force Text format first so Excel does not
# --- auto-convert them into real dates, then strip the temporary format so
# --- the cell keeps the workbook's default (General) style, same as before.
$ws.Range("E2:F2").NumberFormat = "@"
$ws.Range("E2").Value = "10.10.1988"
$ws.Range("F2").Value = "10.10.2015"
$ws.Range("E2:F2").ClearFormats()

# --- Numeric fields on row 2 ---
$ws.Range("C2").Value = 8245615232
$ws.Range("J2").Value = 5
$ws.Range("Q2").Value = 352001
$ws.Range("R2").Value = 5245252
$ws.Range("S2").Value = 415245263
